$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous used range (A1:BA24) contents; ClearContents preserves cell styles (e.g. date format),
# Clear() on the columns/rows being fully removed (column BA, rows 23:24) strips the now-unwanted styles too.
$ws.Range("A1:BA24").ClearContents()
$ws.Range("BA1:BA24").Clear()
$ws.Range("A23:BA24").Clear()

# Row 1: header date serials (B1:AZ1)
$row1 = New-Object 'object[,]' 1,51
$row1[0,0] = 39583
$row1[0,1] = 39765
$row1[0,2] = 39948
$row1[0,3] = 40130
$row1[0,4] = 40310
$row1[0,5] = 40494
$row1[0,6] = 40676
$row1[0,7] = 40862
$row1[0,8] = 41044
$row1[0,9] = 41228
$row1[0,10] = 41409
$row1[0,11] = 41592
$row1[0,12] = 41774
$row1[0,13] = 41957
$row1[0,14] = 42137
$row1[0,15] = 42321
$row1[0,16] = 42503
$row1[0,17] = 42689
$row1[0,18] = 42867
$row1[0,19] = 43053
$row1[0,20] = 43145
$row1[0,21] = 43235
$row1[0,22] = 43326
$row1[0,23] = 43418
$row1[0,24] = 43510
$row1[0,25] = 43600
$row1[0,26] = 43691
$row1[0,27] = 43783
$row1[0,28] = 43875
$row1[0,29] = 43966
$row1[0,30] = 44068
$row1[0,31] = 44159
$row1[0,32] = 44251
$row1[0,33] = 44341
$row1[0,34] = 44432
$row1[0,35] = 44525
$row1[0,36] = 44617
$row1[0,37] = 44706
$row1[0,38] = 44798
$row1[0,39] = 44890
$row1[0,40] = 44981
$row1[0,41] = 45071
$row1[0,42] = 45163
$row1[0,43] = 45254
$row1[0,44] = 45345
$row1[0,45] = 45436
$row1[0,46] = 45534
$row1[0,47] = 45618
$row1[0,48] = 45713
$row1[0,49] = 45800
$row1[0,50] = 45891
$ws.Range("B1:AZ1").Value = $row1

# Column A: row date serials (A2:A22)
$colA = New-Object 'object[,]' 21,1
$colA[0,0] = 39813
$colA[1,0] = 40178
$colA[2,0] = 40543
$colA[3,0] = 40908
$colA[4,0] = 41274
$colA[5,0] = 41639
$colA[6,0] = 42004
$colA[7,0] = 42369
$colA[8,0] = 42735
$colA[9,0] = 43100
$colA[10,0] = 43465
$colA[11,0] = 43830
$colA[12,0] = 44196
$colA[13,0] = 44561
$colA[14,0] = 44926
$colA[15,0] = 45291
$colA[16,0] = 45657
$colA[17,0] = 46022
$colA[18,0] = 46387
$colA[19,0] = 46752
$colA[20,0] = 47118
$ws.Range("A2:A22").Value = $colA

# Forecast data rows (each row is a contiguous range ending at column AZ)
$r3 = New-Object 'object[,]' 1,48
$r3[0,0] = -4.715480642250625
$r3[0,1] = -4.715480642250625
$r3[0,2] = -4.715480642250625
$r3[0,3] = -4.715480642250625
$r3[0,4] = -4.715480642250625
$r3[0,5] = -4.715480642250625
$r3[0,6] = -4.715480642250625
$r3[0,7] = -4.715480642250625
$r3[0,8] = -4.715480642250625
$r3[0,9] = -4.715480642250625
$r3[0,10] = -4.715480642250625
$r3[0,11] = -4.715480642250625
$r3[0,12] = -4.715480642250625
$r3[0,13] = -4.715480642250625
$r3[0,14] = -4.715480642250625
$r3[0,15] = -4.715480642250625
$r3[0,16] = -4.715480642250625
$r3[0,17] = -4.715480642250625
$r3[0,18] = -4.715480642250625
$r3[0,19] = -4.715480642250625
$r3[0,20] = -4.715480642250625
$r3[0,21] = -4.715480642250625
$r3[0,22] = -4.715480642250625
$r3[0,23] = -4.715480642250625
$r3[0,24] = -4.715480642250625
$r3[0,25] = -4.715480642250625
$r3[0,26] = -4.715480642250625
$r3[0,27] = -4.715480642250625
$r3[0,28] = -4.715480642250625
$r3[0,29] = -4.715480642250625
$r3[0,30] = -4.715480642250625
$r3[0,31] = -4.715480642250625
$r3[0,32] = -4.715480642250625
$r3[0,33] = -4.715480642250625
$r3[0,34] = -4.715480642250625
$r3[0,35] = -4.715480642250625
$r3[0,36] = -4.715480642250625
$r3[0,37] = -4.715480642250625
$r3[0,38] = -4.715480642250625
$r3[0,39] = -4.715480642250625
$r3[0,40] = -4.715480642250625
$r3[0,41] = -4.715480642250625
$r3[0,42] = -4.715480642250625
$r3[0,43] = -4.715480642250625
$r3[0,44] = -4.715480642250625
$r3[0,45] = -4.715480642250625
$r3[0,46] = -4.715480642250625
$r3[0,47] = -4.715480642250625
$ws.Range("E3:AZ3").Value = $r3

$r4 = New-Object 'object[,]' 1,46
$r4[0,0] = 6.130685532900881
$r4[0,1] = 6.130685532900881
$r4[0,2] = 6.130685532900881
$r4[0,3] = 6.130685532900881
$r4[0,4] = 6.130685532900904
$r4[0,5] = 6.130685532900904
$r4[0,6] = 6.130685532900904
$r4[0,7] = 6.130685532900904
$r4[0,8] = 6.130685532900904
$r4[0,9] = 6.130685532900904
$r4[0,10] = 6.130685532900904
$r4[0,11] = 6.130685532900904
$r4[0,12] = 6.130685532900904
$r4[0,13] = 6.130685532900904
$r4[0,14] = 6.130685532900904
$r4[0,15] = 6.130685532900904
$r4[0,16] = 6.130685532900904
$r4[0,17] = 6.130685532900904
$r4[0,18] = 6.130685532900904
$r4[0,19] = 6.130685532900904
$r4[0,20] = 6.130685532900904
$r4[0,21] = 6.130685532900904
$r4[0,22] = 6.130685532900904
$r4[0,23] = 6.130685532900904
$r4[0,24] = 6.130685532900904
$r4[0,25] = 6.130685532900904
$r4[0,26] = 6.130685532900904
$r4[0,27] = 6.130685532900904
$r4[0,28] = 6.130685532900904
$r4[0,29] = 6.130685532900904
$r4[0,30] = 6.130685532900904
$r4[0,31] = 6.130685532900904
$r4[0,32] = 6.130685532900904
$r4[0,33] = 6.130685532900904
$r4[0,34] = 6.130685532900904
$r4[0,35] = 6.130685532900904
$r4[0,36] = 6.130685532900904
$r4[0,37] = 6.130685532900904
$r4[0,38] = 6.130685532900904
$r4[0,39] = 6.130685532900904
$r4[0,40] = 6.130685532900904
$r4[0,41] = 6.130685532900904
$r4[0,42] = 6.130685532900904
$r4[0,43] = 6.130685532900904
$r4[0,44] = 6.130685532900904
$r4[0,45] = 6.130685532900904
$ws.Range("G4:AZ4").Value = $r4

$r5 = New-Object 'object[,]' 1,44
$r5[0,0] = 8.703939237319025
$r5[0,1] = 8.703939237319025
$r5[0,2] = 8.703939237318981
$r5[0,3] = 8.703939237318981
$r5[0,4] = 8.703939237318981
$r5[0,5] = 8.703939237318981
$r5[0,6] = 8.703939237318981
$r5[0,7] = 8.703939237318981
$r5[0,8] = 8.703939237318981
$r5[0,9] = 8.703939237318981
$r5[0,10] = 8.703939237318981
$r5[0,11] = 8.703939237318981
$r5[0,12] = 8.703939237318981
$r5[0,13] = 8.703939237318981
$r5[0,14] = 8.703939237318981
$r5[0,15] = 8.703939237318981
$r5[0,16] = 8.703939237318981
$r5[0,17] = 8.703939237318981
$r5[0,18] = 8.703939237318981
$r5[0,19] = 8.703939237318981
$r5[0,20] = 8.703939237318981
$r5[0,21] = 8.703939237318981
$r5[0,22] = 8.703939237318981
$r5[0,23] = 8.703939237318981
$r5[0,24] = 8.703939237318981
$r5[0,25] = 8.703939237318981
$r5[0,26] = 8.703939237318981
$r5[0,27] = 8.703939237318981
$r5[0,28] = 8.703939237318981
$r5[0,29] = 8.703939237318981
$r5[0,30] = 8.703939237318981
$r5[0,31] = 8.703939237318981
$r5[0,32] = 8.703939237318981
$r5[0,33] = 8.703939237318981
$r5[0,34] = 8.703939237318981
$r5[0,35] = 8.703939237318981
$r5[0,36] = 8.703939237318981
$r5[0,37] = 8.703939237318981
$r5[0,38] = 8.703939237318981
$r5[0,39] = 8.703939237318981
$r5[0,40] = 8.703939237318981
$r5[0,41] = 8.703939237318981
$r5[0,42] = 8.703939237318981
$r5[0,43] = 8.703939237318981
$ws.Range("I5:AZ5").Value = $r5

$r6 = New-Object 'object[,]' 1,42
$r6[0,0] = 2.688274587589135
$r6[0,1] = 2.688274587589135
$r6[0,2] = 2.688274587589135
$r6[0,3] = 2.688274587589135
$r6[0,4] = 2.688274587589135
$r6[0,5] = 2.688274587589135
$r6[0,6] = 2.688274587589135
$r6[0,7] = 2.688274587589135
$r6[0,8] = 2.688274587589135
$r6[0,9] = 2.688274587589135
$r6[0,10] = 2.688274587589135
$r6[0,11] = 2.688274587589135
$r6[0,12] = 2.688274587589135
$r6[0,13] = 2.688274587589135
$r6[0,14] = 2.688274587589135
$r6[0,15] = 2.688274587589135
$r6[0,16] = 2.688274587589135
$r6[0,17] = 2.688274587589135
$r6[0,18] = 2.688274587589135
$r6[0,19] = 2.688274587589135
$r6[0,20] = 2.688274587589135
$r6[0,21] = 2.688274587589135
$r6[0,22] = 2.688274587589135
$r6[0,23] = 2.688274587589135
$r6[0,24] = 2.688274587589135
$r6[0,25] = 2.688274587589135
$r6[0,26] = 2.688274587589135
$r6[0,27] = 2.688274587589135
$r6[0,28] = 2.688274587589135
$r6[0,29] = 2.688274587589135
$r6[0,30] = 2.688274587589135
$r6[0,31] = 2.688274587589135
$r6[0,32] = 2.688274587589135
$r6[0,33] = 2.688274587589135
$r6[0,34] = 2.688274587589135
$r6[0,35] = 2.688274587589135
$r6[0,36] = 2.688274587589135
$r6[0,37] = 2.688274587589135
$r6[0,38] = 2.688274587589135
$r6[0,39] = 2.688274587589135
$r6[0,40] = 2.688274587589135
$r6[0,41] = 2.688274587589135
$ws.Range("K6:AZ6").Value = $r6

$r7 = New-Object 'object[,]' 1,42
$r7[0,0] = 1.490702606731831
$r7[0,1] = 0.4454453461194552
$r7[0,2] = 0.9946838291217786
$r7[0,3] = 0.9946838291217786
$r7[0,4] = 0.9946838291217786
$r7[0,5] = 0.9946838291217786
$r7[0,6] = 0.9946838291217786
$r7[0,7] = 0.9946838291217786
$r7[0,8] = 0.9946838291217786
$r7[0,9] = 0.9946838291217786
$r7[0,10] = 0.9946838291217786
$r7[0,11] = 0.9946838291217786
$r7[0,12] = 0.9946838291217786
$r7[0,13] = 0.9946838291217786
$r7[0,14] = 0.9946838291217786
$r7[0,15] = 0.9946838291217786
$r7[0,16] = 0.9946838291217786
$r7[0,17] = 0.9946838291217786
$r7[0,18] = 0.9946838291217786
$r7[0,19] = 0.9946838291217786
$r7[0,20] = 0.9946838291217786
$r7[0,21] = 0.9946838291217786
$r7[0,22] = 0.9946838291217786
$r7[0,23] = 0.9946838291217786
$r7[0,24] = 0.9946838291217786
$r7[0,25] = 0.9946838291217786
$r7[0,26] = 0.9946838291217786
$r7[0,27] = 0.9946838291217786
$r7[0,28] = 0.9946838291217786
$r7[0,29] = 0.9946838291217786
$r7[0,30] = 0.9946838291217786
$r7[0,31] = 0.9946838291217786
$r7[0,32] = 0.9946838291217786
$r7[0,33] = 0.9946838291217786
$r7[0,34] = 0.9946838291217786
$r7[0,35] = 0.9946838291217786
$r7[0,36] = 0.9946838291217786
$r7[0,37] = 0.9946838291217786
$r7[0,38] = 0.9946838291217786
$r7[0,39] = 0.9946838291217786
$r7[0,40] = 0.9946838291217786
$r7[0,41] = 0.9946838291217786
$ws.Range("K7:AZ7").Value = $r7

$r8 = New-Object 'object[,]' 1,42
$r8[0,0] = 2.736923691861937
$r8[0,1] = 2.887668087172179
$r8[0,2] = 3.498411812952029
$r8[0,3] = 5.461771395837989
$r8[0,4] = 5.562499360312567
$r8[0,5] = 5.562499360312567
$r8[0,6] = 5.562499360312567
$r8[0,7] = 5.562499360312567
$r8[0,8] = 5.562499360312567
$r8[0,9] = 5.562499360312567
$r8[0,10] = 5.562499360312567
$r8[0,11] = 5.562499360312567
$r8[0,12] = 5.562499360312567
$r8[0,13] = 5.562499360312567
$r8[0,14] = 5.562499360312567
$r8[0,15] = 5.562499360312567
$r8[0,16] = 5.562499360312567
$r8[0,17] = 5.562499360312567
$r8[0,18] = 5.562499360312567
$r8[0,19] = 5.562499360312567
$r8[0,20] = 5.562499360312567
$r8[0,21] = 5.562499360312567
$r8[0,22] = 5.562499360312567
$r8[0,23] = 5.562499360312567
$r8[0,24] = 5.562499360312567
$r8[0,25] = 5.562499360312567
$r8[0,26] = 5.562499360312567
$r8[0,27] = 5.562499360312567
$r8[0,28] = 5.562499360312567
$r8[0,29] = 5.562499360312567
$r8[0,30] = 5.562499360312567
$r8[0,31] = 5.562499360312567
$r8[0,32] = 5.562499360312567
$r8[0,33] = 5.562499360312567
$r8[0,34] = 5.562499360312567
$r8[0,35] = 5.562499360312567
$r8[0,36] = 5.562499360312567
$r8[0,37] = 5.562499360312567
$r8[0,38] = 5.562499360312567
$r8[0,39] = 5.562499360312567
$r8[0,40] = 5.562499360312567
$r8[0,41] = 5.562499360312567
$ws.Range("K8:AZ8").Value = $r8

$r9 = New-Object 'object[,]' 1,41
$r9[0,0] = 2.665853845175659
$r9[0,1] = 2.796235363530908
$r9[0,2] = 3.667647087004666
$r9[0,3] = 4.543069198269034
$r9[0,4] = 3.857137494611718
$r9[0,5] = 4.195080504802551
$r9[0,6] = 4.195080504802551
$r9[0,7] = 4.195080504802551
$r9[0,8] = 4.195080504802551
$r9[0,9] = 4.195080504802551
$r9[0,10] = 4.195080504802551
$r9[0,11] = 4.195080504802551
$r9[0,12] = 4.195080504802551
$r9[0,13] = 4.195080504802551
$r9[0,14] = 4.195080504802551
$r9[0,15] = 4.195080504802551
$r9[0,16] = 4.195080504802551
$r9[0,17] = 4.195080504802551
$r9[0,18] = 4.195080504802551
$r9[0,19] = 4.195080504802551
$r9[0,20] = 4.195080504802551
$r9[0,21] = 4.195080504802551
$r9[0,22] = 4.195080504802551
$r9[0,23] = 4.195080504802551
$r9[0,24] = 4.195080504802551
$r9[0,25] = 4.195080504802551
$r9[0,26] = 4.195080504802551
$r9[0,27] = 4.195080504802551
$r9[0,28] = 4.195080504802551
$r9[0,29] = 4.195080504802551
$r9[0,30] = 4.195080504802551
$r9[0,31] = 4.195080504802551
$r9[0,32] = 4.195080504802551
$r9[0,33] = 4.195080504802551
$r9[0,34] = 4.195080504802551
$r9[0,35] = 4.195080504802551
$r9[0,36] = 4.195080504802551
$r9[0,37] = 4.195080504802551
$r9[0,38] = 4.195080504802551
$r9[0,39] = 4.195080504802551
$r9[0,40] = 4.195080504802551
$ws.Range("L9:AZ9").Value = $r9

$r10 = New-Object 'object[,]' 1,39
$r10[0,0] = 3.241922060312707
$r10[0,1] = 3.457193745607912
$r10[0,2] = 3.60208822706134
$r10[0,3] = 4.42512316868644
$r10[0,4] = 4.325618632128836
$r10[0,5] = 4.230623896992025
$r10[0,6] = 4.230623896992025
$r10[0,7] = 4.230623896992025
$r10[0,8] = 4.230623896992025
$r10[0,9] = 4.230623896992025
$r10[0,10] = 4.230623896992025
$r10[0,11] = 4.230623896992025
$r10[0,12] = 4.230623896992025
$r10[0,13] = 4.230623896992025
$r10[0,14] = 4.230623896992025
$r10[0,15] = 4.230623896992025
$r10[0,16] = 4.230623896992025
$r10[0,17] = 4.230623896992025
$r10[0,18] = 4.230623896992025
$r10[0,19] = 4.230623896992025
$r10[0,20] = 4.230623896992025
$r10[0,21] = 4.230623896992025
$r10[0,22] = 4.230623896992025
$r10[0,23] = 4.230623896992025
$r10[0,24] = 4.230623896992025
$r10[0,25] = 4.230623896992025
$r10[0,26] = 4.230623896992025
$r10[0,27] = 4.230623896992025
$r10[0,28] = 4.230623896992025
$r10[0,29] = 4.230623896992025
$r10[0,30] = 4.230623896992025
$r10[0,31] = 4.230623896992025
$r10[0,32] = 4.230623896992025
$r10[0,33] = 4.230623896992025
$r10[0,34] = 4.230623896992025
$r10[0,35] = 4.230623896992025
$r10[0,36] = 4.230623896992025
$r10[0,37] = 4.230623896992025
$r10[0,38] = 4.230623896992025
$ws.Range("N10:AZ10").Value = $r10

$r11 = New-Object 'object[,]' 1,37
$r11[0,0] = 3.33541401017492
$r11[0,1] = 3.497684176291571
$r11[0,2] = 3.998755562728684
$r11[0,3] = 4.210645455310114
$r11[0,4] = 4.439146757103352
$r11[0,5] = 4.933871867981643
$r11[0,6] = 4.933871867981643
$r11[0,7] = 4.933871867981643
$r11[0,8] = 4.933871867981643
$r11[0,9] = 4.933871867981643
$r11[0,10] = 4.933871867981643
$r11[0,11] = 4.933871867981643
$r11[0,12] = 4.933871867981643
$r11[0,13] = 4.933871867981643
$r11[0,14] = 4.933871867981643
$r11[0,15] = 4.933871867981643
$r11[0,16] = 4.933871867981643
$r11[0,17] = 4.933871867981643
$r11[0,18] = 4.933871867981643
$r11[0,19] = 4.933871867981643
$r11[0,20] = 4.933871867981643
$r11[0,21] = 4.933871867981643
$r11[0,22] = 4.933871867981643
$r11[0,23] = 4.933871867981643
$r11[0,24] = 4.933871867981643
$r11[0,25] = 4.933871867981643
$r11[0,26] = 4.933871867981643
$r11[0,27] = 4.933871867981643
$r11[0,28] = 4.933871867981643
$r11[0,29] = 4.933871867981643
$r11[0,30] = 4.933871867981643
$r11[0,31] = 4.933871867981643
$r11[0,32] = 4.933871867981643
$r11[0,33] = 4.933871867981643
$r11[0,34] = 4.933871867981643
$r11[0,35] = 4.933871867981643
$r11[0,36] = 4.933871867981643
$ws.Range("P11:AZ11").Value = $r11

$r12 = New-Object 'object[,]' 1,35
$r12[0,0] = 3.527005305694031
$r12[0,1] = 3.584738032309165
$r12[0,2] = 3.861679870292711
$r12[0,3] = 4.757571096183799
$r12[0,4] = 5.787381971961936
$r12[0,5] = 5.723509166364238
$r12[0,6] = 5.222860865675738
$r12[0,7] = 5.456119081407906
$r12[0,8] = 5.456119081407906
$r12[0,9] = 5.456119081407906
$r12[0,10] = 5.456119081407906
$r12[0,11] = 5.456119081407906
$r12[0,12] = 5.456119081407906
$r12[0,13] = 5.456119081407906
$r12[0,14] = 5.456119081407906
$r12[0,15] = 5.456119081407906
$r12[0,16] = 5.456119081407906
$r12[0,17] = 5.456119081407906
$r12[0,18] = 5.456119081407906
$r12[0,19] = 5.456119081407906
$r12[0,20] = 5.456119081407906
$r12[0,21] = 5.456119081407906
$r12[0,22] = 5.456119081407906
$r12[0,23] = 5.456119081407906
$r12[0,24] = 5.456119081407906
$r12[0,25] = 5.456119081407906
$r12[0,26] = 5.456119081407906
$r12[0,27] = 5.456119081407906
$r12[0,28] = 5.456119081407906
$r12[0,29] = 5.456119081407906
$r12[0,30] = 5.456119081407906
$r12[0,31] = 5.456119081407906
$r12[0,32] = 5.456119081407906
$r12[0,33] = 5.456119081407906
$r12[0,34] = 5.456119081407906
$ws.Range("R12:AZ12").Value = $r12

$r13 = New-Object 'object[,]' 1,33
$r13[0,0] = 3.594510929080963
$r13[0,1] = 3.734834974488588
$r13[0,2] = 4.118094101621717
$r13[0,3] = 4.058053416301188
$r13[0,4] = 2.939060555390971
$r13[0,5] = 4.569144243718659
$r13[0,6] = 3.466212706516147
$r13[0,7] = 4.674926984813466
$r13[0,8] = 3.038115835571786
$r13[0,9] = 3.346849276607955
$r13[0,10] = 3.346849276607955
$r13[0,11] = 3.346849276607955
$r13[0,12] = 3.346849276607955
$r13[0,13] = 3.346849276607955
$r13[0,14] = 3.346849276607955
$r13[0,15] = 3.346849276607955
$r13[0,16] = 3.346849276607955
$r13[0,17] = 3.346849276607955
$r13[0,18] = 3.346849276607955
$r13[0,19] = 3.346849276607955
$r13[0,20] = 3.346849276607955
$r13[0,21] = 3.346849276607955
$r13[0,22] = 3.346849276607955
$r13[0,23] = 3.346849276607955
$r13[0,24] = 3.346849276607955
$r13[0,25] = 3.346849276607955
$r13[0,26] = 3.346849276607955
$r13[0,27] = 3.346849276607955
$r13[0,28] = 3.346849276607955
$r13[0,29] = 3.346849276607955
$r13[0,30] = 3.346849276607955
$r13[0,31] = 3.346849276607955
$r13[0,32] = 3.346849276607955
$ws.Range("T13:AZ13").Value = $r13

$r14 = New-Object 'object[,]' 1,30
$r14[0,0] = 3.849279581596932
$r14[0,1] = 3.673587953629398
$r14[0,2] = 3.98814479861731
$r14[0,3] = 3.597750881470851
$r14[0,4] = 4.726969153629335
$r14[0,5] = 1.154972712087221
$r14[0,6] = 3.124801698476176
$r14[0,7] = 3.370990011762443
$r14[0,8] = 0.4167846160013644
$r14[0,9] = -9.2489161297999
$r14[0,10] = -9.2489161297999
$r14[0,11] = -9.2489161297999
$r14[0,12] = -9.2489161297999
$r14[0,13] = -9.2489161297999
$r14[0,14] = -9.2489161297999
$r14[0,15] = -9.2489161297999
$r14[0,16] = -9.2489161297999
$r14[0,17] = -9.2489161297999
$r14[0,18] = -9.2489161297999
$r14[0,19] = -9.2489161297999
$r14[0,20] = -9.2489161297999
$r14[0,21] = -9.2489161297999
$r14[0,22] = -9.2489161297999
$r14[0,23] = -9.2489161297999
$r14[0,24] = -9.2489161297999
$r14[0,25] = -9.2489161297999
$r14[0,26] = -9.2489161297999
$r14[0,27] = -9.2489161297999
$r14[0,28] = -9.2489161297999
$r14[0,29] = -9.2489161297999
$ws.Range("W14:AZ14").Value = $r14

$r15 = New-Object 'object[,]' 1,26
$r15[0,0] = 3.97322668464466
$r15[0,1] = 3.352837059134517
$r15[0,2] = 3.754874824891763
$r15[0,3] = 3.843786543692795
$r15[0,4] = 1.125570778878981
$r15[0,5] = -10.65745199005891
$r15[0,6] = -4.101394328717845
$r15[0,7] = -1.861534891151506
$r15[0,8] = -1.488064879190421
$r15[0,9] = -1.287084480507283
$r15[0,10] = -1.287084480507283
$r15[0,11] = -1.287084480507283
$r15[0,12] = -1.287084480507283
$r15[0,13] = -1.287084480507283
$r15[0,14] = -1.287084480507283
$r15[0,15] = -1.287084480507283
$r15[0,16] = -1.287084480507283
$r15[0,17] = -1.287084480507283
$r15[0,18] = -1.287084480507283
$r15[0,19] = -1.287084480507283
$r15[0,20] = -1.287084480507283
$r15[0,21] = -1.287084480507283
$r15[0,22] = -1.287084480507283
$r15[0,23] = -1.287084480507283
$r15[0,24] = -1.287084480507283
$r15[0,25] = -1.287084480507283
$ws.Range("AA15:AZ15").Value = $r15

$r16 = New-Object 'object[,]' 1,22
$r16[0,0] = 3.197223976813479
$r16[0,1] = 2.090835525509216
$r16[0,2] = 2.211767994230307
$r16[0,3] = 3.108682697521514
$r16[0,4] = 3.63609986063671
$r16[0,5] = 4.124305474197043
$r16[0,6] = 3.801772939051373
$r16[0,7] = 1.920033066224791
$r16[0,8] = 1.712986619197032
$r16[0,9] = 1.494343500592232
$r16[0,10] = 1.494343500592232
$r16[0,11] = 1.494343500592232
$r16[0,12] = 1.494343500592232
$r16[0,13] = 1.494343500592232
$r16[0,14] = 1.494343500592232
$r16[0,15] = 1.494343500592232
$r16[0,16] = 1.494343500592232
$r16[0,17] = 1.494343500592232
$r16[0,18] = 1.494343500592232
$r16[0,19] = 1.494343500592232
$r16[0,20] = 1.494343500592232
$r16[0,21] = 1.494343500592232
$ws.Range("AE16:AZ16").Value = $r16

$r17 = New-Object 'object[,]' 1,19
$r17[0,0] = 3.042266670461236
$r17[0,1] = 3.386238465899294
$r17[0,2] = 3.537188017584847
$r17[0,3] = 3.4135518151031
$r17[0,4] = 2.774323849124349
$r17[0,5] = 2.431967849366434
$r17[0,6] = 1.882230726672129
$r17[0,7] = 0.9099262091262217
$r17[0,8] = -0.07673633990846751
$r17[0,9] = -0.5717743519535134
$r17[0,10] = -0.6982718287330991
$r17[0,11] = -0.6982718287330991
$r17[0,12] = -0.6982718287330991
$r17[0,13] = -0.6982718287330991
$r17[0,14] = -0.6982718287330991
$r17[0,15] = -0.6982718287330991
$r17[0,16] = -0.6982718287330991
$r17[0,17] = -0.6982718287330991
$r17[0,18] = -0.6982718287330991
$ws.Range("AH17:AZ17").Value = $r17

$r18 = New-Object 'object[,]' 1,15
$r18[0,0] = 3.109945199678088
$r18[0,1] = 2.846831818037354
$r18[0,2] = 2.624147399924981
$r18[0,3] = 2.509020922481398
$r18[0,4] = 2.180844122535164
$r18[0,5] = 1.732880403074311
$r18[0,6] = 1.415512869596025
$r18[0,7] = 0.1232424362653362
$r18[0,8] = -0.5700058398449448
$r18[0,9] = -0.2867681914691111
$r18[0,10] = -0.4137309550271362
$r18[0,11] = -0.4137309550271362
$r18[0,12] = -0.4137309550271362
$r18[0,13] = -0.4137309550271362
$r18[0,14] = -0.4137309550271362
$ws.Range("AL18:AZ18").Value = $r18

$r19 = New-Object 'object[,]' 1,11
$r19[0,0] = 2.277761159626501
$r19[0,1] = 2.080657666875685
$r19[0,2] = 1.947454963355622
$r19[0,3] = 1.970932490584887
$r19[0,4] = 1.804123797928292
$r19[0,5] = 2.013081730696564
$r19[0,6] = 1.790585695398428
$r19[0,7] = 1.743978804508384
$r19[0,8] = 0.946259771301472
$r19[0,9] = 0.4432539413513181
$r19[0,10] = 0.2267356977060819
$ws.Range("AP19:AZ19").Value = $r19

$r20 = New-Object 'object[,]' 1,7
$r20[0,0] = 1.825245137774356
$r20[0,1] = 1.88619841018669
$r20[0,2] = 1.881221885219464
$r20[0,3] = 1.964863231054204
$r20[0,4] = 1.770613036357038
$r20[0,5] = 1.364481450639365
$r20[0,6] = 0.9098136509666066
$ws.Range("AT20:AZ20").Value = $r20

$r21 = New-Object 'object[,]' 1,3
$r21[0,0] = 1.898066514521801
$r21[0,1] = 1.754988417644232
$r21[0,2] = 1.645330300307712
$ws.Range("AX21:AZ21").Value = $r21
